# "Sensoren auslesen, Zeitaufwand minimiert"
# Add a new "Sensoren Optimieren" column (I) to the top summary table and
# two new detail columns (F/H) to the lower "Arduino Optimierungen" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lower table: new "ADS Bypass" (F) / "MPU Lesen" (H) detail columns ---
$ws.Range("F13").Value = "ADS Bypass"
$ws.Range("F14").Value = "1,6ms"
$ws.Range("F15").Value = "0,5ms"
$ws.Range("F16").Value = "1,0ms"
$ws.Range("F17").Value = "0ms"
$ws.Range("F19").Value = "Pressure weg"

$ws.Range("H13").Value = "MPU Lesen"
$ws.Range("H14").Value = "3ms"
$ws.Range("H15").Value = "0,5ms"
$ws.Range("H16").Value = "2,5ms"
$ws.Range("H17").Value = "0ms"

# --- Upper table: new "Sensoren Optimieren" column (I) ---
$ws.Range("I2").Value = "Sensoren Optimieren"
$ws.Range("I7").Value = "'---"
$ws.Range("I4").Value = "ca. 6ms"
$ws.Range("I5").Value = "2,1ms"
$ws.Range("I6").Value = "3ms"

# --- View state: move selection to J24 ---
[void]$ws.Range("J24").Select()
